# Proof reading of chapter Introduction and OpenCL
# Add a new log entry (date + activity) to the protocol sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row: date serial 41385 (2013-04-21) in A71, text in B71.
$ws.Cells.Item(71, 1).Value = 41385
$ws.Cells.Item(71, 2).Value = "proof-read introduction and OpenCL chapter"

# Move the active selection to B73, matching where the user ended up.
$ws.Range("B73").Select()
